# Update "want to go" counts (column F) on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, matching the published site refresh.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes    = $wb.Worksheets.Item(4)   # 全部类型

# Row -> new F value for the "展览" sheet (rows are 1-based, header = row 1)
$exhibitionUpdates = @{
    4  = 81
    5  = 18
    6  = 543
    7  = 1669
    8  = 19
    10 = 27
    11 = 1615
    13 = 66
    14 = 399
    15 = 261
    16 = 194
    19 = 30
    21 = 207
    22 = 290
    23 = 158
    24 = 223
    25 = 225
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet (one extra row vs. "展览"
# because it also includes the single "演出" entry)
$allTypesUpdates = @{
    4  = 81
    5  = 18
    6  = 543
    7  = 1669
    9  = 19
    11 = 27
    12 = 1615
    14 = 66
    15 = 399
    16 = 261
    17 = 194
    20 = 30
    22 = 207
    23 = 290
    24 = 158
    25 = 223
    26 = 225
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
